$p = $ppt.ActivePresentation

# 1. Change the table style on slide 6's table shape.
$s = $p.Slides.Item(6)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{7491DB22-4DF4-4D2F-8550-F5205CEBE5C4}")

# 2. Swap theme1.xml <-> theme2.xml content (slide master theme <-> notes master theme).
$p.Designs.Item(1).SlideMaster.Theme.Name
